$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/add Experience (C3) and Salary (D3) for row 3
$ws.Range("D3").Value = 800000

# Row 4: add Age (B4)
$ws.Range("B4").Value = 33

# Row 6: Experience changes from 0 to 1
$ws.Range("C6").Value = 1

# Row 7: add Experience (C7)
$ws.Range("C7").Value = 1

# Row 8: add Experience (C8) and Salary (D8)
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0

# Row 9 (Zainab): add Age (B9), Experience (C9), Salary (D9)
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

# Delete old rows 10 and 11 entirely (shrinks dimension to A1:D9)
$ws.Rows("10:11").Delete()

# Update selection to match target state
$ws.Range("D8").Select()
